# Apply the dated replacements described by the diff.
$d = $word.ActiveDocument

# Mapping of old text -> new text (order matches the diff).
$replacements = [ordered]@{
    "2025-04-23 Wednesday" = "2025-04-24 Thursday"
    "61×30="                = "37×40="
    "98×60="                = "87×13="
    "32×63="                = "91×93="
    "85×96="                = "27×58="
    "88×23="                = "17×79="
    "15×33="                = "89×75="
    "84×74="                = "94×18="
    "20×19="                = "97×53="
    "18×17="                = "44×61="
    "76×19="                = "53×49="
    "43×66="                = "64×68="
    "23×25="                = "31×92="
    "64×49="                = "76×76="
    "27×93="                = "85×36="
    "59×41="                = "80×31="
    "85×87="                = "29×73="
    "40×47="                = "21×92="
    "49×35="                = "55×95="
    "31×27="                = "63×93="
    "42×78="                = "81×54="
    "65×37="                = "11×24="
    "75×79="                = "48×80="
    "81×73="                = "76×85="
    "75×39="                = "20×13="
    "39×24="                = "65×13="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
